$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 428:429, pushing the existing rows 428-446 down to 430-448
$ws.Rows("428:429").Insert()

# Row 428 - new weekly data (Primera)
$ws.Range("A428").Value = 8
$ws.Range("B428").Value = "Terminal La Palmera de La Serena"
$ws.Range("C428").Value = "Coquimbo"
$ws.Range("D428").Value = 44753
$ws.Range("E428").Value = 4
$ws.Range("F428").Value = 100112043
$ws.Range("G428").Value = "Pepino dulce"
$ws.Range("H428").Value = "Cultivar IV Región"
$ws.Range("I428").Value = "Primera"
$ws.Range("J428").Value = 520
$ws.Range("K428").Value = 13000
$ws.Range("L428").Value = 14000
$ws.Range("M428").Value = 13500
$ws.Range("N428").Value = "$/bandeja 18 kilos"
$ws.Range("O428").Value = "Provincia de Limarí"
$ws.Range("P428").Value = 750
$ws.Range("Q428").Value = 18
$ws.Range("R428").Value = "Hortaliza"

# Row 429 - new weekly data (Segunda)
$ws.Range("A429").Value = 8
$ws.Range("B429").Value = "Terminal La Palmera de La Serena"
$ws.Range("C429").Value = "Coquimbo"
$ws.Range("D429").Value = 44753
$ws.Range("E429").Value = 4
$ws.Range("F429").Value = 100112043
$ws.Range("G429").Value = "Pepino dulce"
$ws.Range("H429").Value = "Cultivar IV Región"
$ws.Range("I429").Value = "Segunda"
$ws.Range("J429").Value = 360
$ws.Range("K429").Value = 10000
$ws.Range("L429").Value = 11000
$ws.Range("M429").Value = 10500
$ws.Range("N429").Value = "$/bandeja 18 kilos"
$ws.Range("O429").Value = "Provincia de Limarí"
$ws.Range("P429").Value = 583
$ws.Range("Q429").Value = 18
$ws.Range("R429").Value = "Hortaliza"

# Apply the date number format (style index 2 / "YYYY-MM-DD HH:MM:SS") used throughout column D
$ws.Range("D428").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D429").NumberFormat = "YYYY-MM-DD HH:MM:SS"
